$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data values. Cells are stored as text (General format)
# in the original sheet, so we force text via NumberFormat "@" before writing
# to avoid Excel auto-converting numeric-looking strings (e.g. "308.52", "19",
# "0.07621") into numeric cells, then reset the style back to Normal so no
# stray cell-style index is left behind (matches original unstyled cells).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '308.52'
$c.Style = "Normal"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '40.81'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '0.12%'
$c.Style = "Normal"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.118'
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '1.25%'
$c.Style = "Normal"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.07621'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '0.05%'
$c.Style = "Normal"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B6")
$c.NumberFormat = "@"
$c.Value = 'GateToken'
$c.Style = "Normal"
$c = $ws.Range("C6")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '4.286'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '0.97%'
$c.Style = "Normal"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = 'FTXToken'
$c.Style = "Normal"
$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.615'
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '1.43%'
$c.Style = "Normal"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = 'BTSEToken'
$c.Style = "Normal"
$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '2.456'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '0.66%'
$c.Style = "Normal"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = 'MXToken'
$c.Style = "Normal"
$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.9075'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '0.11%'
$c.Style = "Normal"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = 'LiechtensteinCryptoassetsExchange'
$c.Style = "Normal"
$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.1285'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '29.94%'
$c.Style = "Normal"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = 'WazirX'
$c.Style = "Normal"
$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.1801'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '1.97%'
$c.Style = "Normal"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = 'MandalaExchangeToken'
$c.Style = "Normal"
$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.09067'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '-0.24%'
$c.Style = "Normal"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = 'BitrueCoin'
$c.Style = "Normal"
$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.04329'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '-1.23%'
$c.Style = "Normal"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = 'BitMartToken'
$c.Style = "Normal"
$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.1044'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '-0.86%'
$c.Style = "Normal"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = 'BitForexToken'
$c.Style = "Normal"
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.001255'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '0.35%'
$c.Style = "Normal"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = 'TigerCash'
$c.Style = "Normal"
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.005746'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '-0.92%'
$c.Style = "Normal"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = 'LEO'
$c.Style = "Normal"
$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '3.345'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '-0.63%'
$c.Style = "Normal"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '0.51%'
$c.Style = "Normal"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.951'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '1.63%'
$c.Style = "Normal"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '4.03%'
$c.Style = "Normal"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '-4.82%'
$c.Style = "Normal"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.04043'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '-2.71%'
$c.Style = "Normal"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.001273'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '4.92%'
$c.Style = "Normal"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.004047'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '0.97%'
$c.Style = "Normal"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.0001272'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '-2.15%'
$c.Style = "Normal"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '24.69%'
$c.Style = "Normal"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02422'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '0.14%'
$c.Style = "Normal"
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05213'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '1.35%'
$c.Style = "Normal"
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.007835'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '0.03%'
$c.Style = "Normal"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.1299'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '-0.65%'
$c.Style = "Normal"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.006805'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-3.69%'
$c.Style = "Normal"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.001934'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '-0.77%'
$c.Style = "Normal"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.007365'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '-8.35%'
$c.Style = "Normal"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.3349'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '9.74%'
$c.Style = "Normal"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00006892'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '8.05%'
$c.Style = "Normal"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.1103'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '1,829.62%'
$c.Style = "Normal"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0002004'
$c.Style = "Normal"
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '19'
$c.Style = "Normal"
